$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Row 2 headers: "unnamed: 1_level_1" and "unnamed: 5_level_1" were placeholder
# labels produced by a multi-index export; correct them to "total" to match
# the other merged "total" column labels.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
